$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "81.395.30"
$ws.Range("E2").Value = "  +5.28%  "
$ws.Range("D3").Value = "3.188.46"
$ws.Range("E3").Value = "  +1.49%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'210.00"
$ws.Range("E5").Value = "  +3.73%  "
$ws.Range("D6").Value = "'635.90"
$ws.Range("E6").Value = "  +0.93%  "
$ws.Range("D7").Value = "'0.291"
$ws.Range("E7").Value = "  +27.97%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.593"
$ws.Range("E9").Value = "  +3.51%  "
$ws.Range("D10").Value = "3.186.61"
$ws.Range("E10").Value = "  +1.50%  "
$ws.Range("D11").Value = "'0.593"
$ws.Range("E11").Value = "  +11.33%  "
$ws.Range("E12").Value = "  +18.98%  "
$ws.Range("E13").Value = "  +2.28%  "
$ws.Range("E14").Value = "  -0.67%  "
$ws.Range("D15").Value = "3.773.31"
$ws.Range("E15").Value = "  +1.49%  "
$ws.Range("E16").Value = "  +5.42%  "
$ws.Range("D17").Value = "81.263.76"
$ws.Range("E17").Value = "  +5.21%  "
$ws.Range("D18").Value = "3.176.76"
$ws.Range("E18").Value = "  +1.49%  "
$ws.Range("D19").Value = "'3.26"
$ws.Range("E19").Value = "  +14.67%  "
$ws.Range("D20").Value = "'14.38"
$ws.Range("E20").Value = "  +3.89%  "
$ws.Range("D21").Value = "'9.28"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("D22").Value = "'442.04"
$ws.Range("E22").Value = "  +3.11%  "
$ws.Range("E23").Value = "  +10.11%  "
$ws.Range("D24").Value = "'7.10"
$ws.Range("E24").Value = "  +5.26%  "
$ws.Range("D25").Value = "'5.08"
$ws.Range("E25").Value = "  +9.27%  "
$ws.Range("D26").Value = "'11.33"
$ws.Range("E26").Value = "  +6.37%  "
$ws.Range("D27").Value = "3.353.24"
$ws.Range("E27").Value = "  +1.55%  "
$ws.Range("D28").Value = "'77.14"
$ws.Range("E28").Value = "  +2.21%  "
$ws.Range("D29").Value = "'0.0000128"
$ws.Range("E29").Value = "  +12.26%  "
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("E31").Value = "  +5.67%  "
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("D33").Value = "'574.91"
$ws.Range("E33").Value = "  +10.63%  "
$ws.Range("E34").Value = "  +2.47%  "
$ws.Range("E35").Value = "  +34.02%  "
$ws.Range("D36").Value = "'2.05"
$ws.Range("E36").Value = "  +4.56%  "
$ws.Range("E37").Value = "  +12.36%  "
$ws.Range("D38").Value = "'23.23"
$ws.Range("E38").Value = "  +4.75%  "
$ws.Range("D39").Value = "'0.999"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").Value = "'0.416"
$ws.Range("E40").Value = "  +5.96%  "
$ws.Range("D41").Value = "'3.13"
$ws.Range("E41").Value = "  +24.55%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").Value = "'5.99"
$ws.Range("E42").Value = "  +11.54%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'2.06"
$ws.Range("E43").Value = "  +17.97%  "
$ws.Range("D44").Value = "'20.80"
$ws.Range("E44").Value = "  +3.69%  "
$ws.Range("D45").Value = "'160.23"
$ws.Range("E45").Value = "  -2.10%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").Value = "'189.66"
$ws.Range("E47").Value = "  -2.76%  "
$ws.Range("D48").Value = "'45.30"
$ws.Range("E48").Value = "  +6.29%  "
$ws.Range("D49").Value = "'1.36"
$ws.Range("E49").Value = "  +5.89%  "
$ws.Range("D50").Value = "'0.785"
$ws.Range("E50").Value = "  -1.65%  "
$ws.Range("D51").Value = "'4.32"
$ws.Range("E51").Value = "  +6.01%  "
